# Add a new "Process_WorkItemsType" entry to the Constants table on the
# "Constants" sheet, storing the work item type that the BOT should
# process (commit: "Added the logic to store the work item type in the
# config file").
#
# The Constants sheet holds a Name / Value / Description table. The new
# entry needs to be inserted right after the existing
# "ACMESystem1_CredentialName" row (row 13) and before the
# "SMTP_OrchestratorAsset" row, keeping the same blank-row-as-separator
# layout that's used throughout the sheet:
#   ...
#   13  ACMESystem1_CredentialName | ACMESystem1_Credentials | ...
#   14  <blank separator row>
#   15  Process_WorkItemsType      | WI2                     | Type of the work items...
#   16  <blank separator row>
#   17  SMTP_OrchestratorAsset     | SMTP_OutlookCredentials | ...
#   ...

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Insert two new rows above the existing row 14 (the blank separator that
# used to sit right below "ACMESystem1_CredentialName"). This pushes the
# old row 14 onward down by two rows, reproducing it as the new blank
# separator row 16, and frees up rows 14-15 for the new entry + its own
# separator.
$ws.Rows.Item(14).Resize(2).Insert()

# The newly inserted rows don't keep the 14.25pt custom row height used
# elsewhere on this sheet, so restore it explicitly.
$ws.Rows.Item(14).Resize(2).RowHeight = 14.25

# Populate the new row (row 15) with the Name / Value / Description for
# the new config entry. Description/Name/Value order mirrors how the
# original authoring tool appended these strings to the shared string
# table.
$ws.Cells.Item(15, 3).Value() = "Type of the work items to be processed by the BOT (E.g: WI2, WI3 e.t.c)"
$ws.Cells.Item(15, 1).Value() = "Process_WorkItemsType"
$ws.Cells.Item(15, 2).Value() = "WI2"

# Reflect the author's final cursor position on the sheet.
[void]$ws.Activate()
[void]$ws.Range("B16").Select()
